$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting existing N:P -> O:Q
$ws.Columns("N:N").Insert()

# The newly inserted column keeps column M's custom width (11 chars)
$ws.Columns("N:N").ColumnWidth = 10.15

# Make "Repayment schedule" the active sheet/tab with R9 selected
$ws.Activate()
$ws.Range("R9").Select() | Out-Null
